$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-12-17T07:01:46.229230+00:00"
$ws.Range("K3").Value = "2025-12-17T07:01:46.229264+00:00"
$ws.Range("K4").Value = "2025-12-17T07:01:48.767228+00:00"
$ws.Range("K5").Value = "2025-12-17T07:01:48.767258+00:00"
$ws.Range("K6").Value = "2025-12-17T07:01:48.767276+00:00"
$ws.Range("K7").Value = "2025-12-17T07:01:51.578433+00:00"
$ws.Range("K8").Value = "2025-12-17T07:01:53.890470+00:00"
$ws.Range("K9").Value = "2025-12-17T07:01:56.656200+00:00"
$ws.Range("K10").Value = "2025-12-17T07:01:56.656229+00:00"
$ws.Range("K11").Value = "2025-12-17T07:01:59.407506+00:00"
$ws.Range("K12").Value = "2025-12-17T07:02:04.038315+00:00"
$ws.Range("K13").Value = "2025-12-17T07:02:04.038344+00:00"
$ws.Range("K14").Value = "2025-12-17T07:02:06.310210+00:00"
$ws.Range("K15").Value = "2025-12-17T07:02:08.629457+00:00"
$ws.Range("K16").Value = "2025-12-17T07:02:10.926875+00:00"
$ws.Range("K17").Value = "2025-12-17T07:02:14.000297+00:00"
$ws.Range("K18").Value = "2025-12-17T07:02:14.000335+00:00"
$ws.Range("K19").Value = "2025-12-17T07:02:14.000364+00:00"
$ws.Range("K20").Value = "2025-12-17T07:02:14.000388+00:00"
$ws.Range("K21").Value = "2025-12-17T07:02:16.303479+00:00"
$ws.Range("K22").Value = "2025-12-17T07:02:16.303508+00:00"
$ws.Range("K23").Value = "2025-12-17T07:02:18.558919+00:00"
$ws.Range("K24").Value = "2025-12-17T07:02:18.558947+00:00"
$ws.Range("K25").Value = "2025-12-17T07:02:18.558967+00:00"
$ws.Range("K26").Value = "2025-12-17T07:02:20.858730+00:00"
$ws.Range("K27").Value = "2025-12-17T07:02:23.088363+00:00"
$ws.Range("K28").Value = "2025-12-17T07:02:23.088400+00:00"
$ws.Range("K29").Value = "2025-12-17T07:02:23.088419+00:00"
$ws.Range("K30").Value = "2025-12-17T07:02:25.261667+00:00"
$ws.Range("K31").Value = "2025-12-17T07:02:27.556403+00:00"
$ws.Range("K32").Value = "2025-12-17T07:02:27.556432+00:00"
$ws.Range("K33").Value = "2025-12-17T07:02:32.180705+00:00"
$ws.Range("K34").Value = "2025-12-17T07:02:32.180734+00:00"
$ws.Range("K35").Value = "2025-12-17T07:02:34.626910+00:00"
$ws.Range("K36").Value = "2025-12-17T07:02:34.626939+00:00"
